$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for rows 2-8 from serial date 45175 to 45183
$ws.Range("C2:C8").Value = 45183
